$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-14 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-15 Sunday", 2) | Out-Null
$d.Content.Find.Execute("27÷2=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=44, 1", 2) | Out-Null
$d.Content.Find.Execute("28÷2=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "14÷8=1, 6", 2) | Out-Null
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "17÷7=2, 3", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=11, 3", 2) | Out-Null
$d.Content.Find.Execute("31÷4=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2) | Out-Null
$d.Content.Find.Execute("36÷9=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=12, 0", 2) | Out-Null
$d.Content.Find.Execute("31÷7=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "46÷5=9, 1", 2) | Out-Null
$d.Content.Find.Execute("29÷2=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷4=16, 3", 2) | Out-Null
$d.Content.Find.Execute("72÷7=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=11, 4", 2) | Out-Null
$d.Content.Find.Execute("83÷4=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=4, 2", 2) | Out-Null
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷2=31, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷8=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=15, 2", 2) | Out-Null
$d.Content.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=15, 1", 2) | Out-Null
$d.Content.Find.Execute("34÷9=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=19, 1", 2) | Out-Null
$d.Content.Find.Execute("33÷8=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "92÷3=30, 2", 2) | Out-Null
$d.Content.Find.Execute("40÷7=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=16, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷4=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=20, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷4=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=11, 4", 2) | Out-Null
$d.Content.Find.Execute("96÷9=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=44, 1", 2) | Out-Null
$d.Content.Find.Execute("21÷5=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷7=4, 6", 2) | Out-Null
$d.Content.Find.Execute("66÷3=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷4=22, 2", 2) | Out-Null
$d.Content.Find.Execute("75÷5=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷3=31, 0", 2) | Out-Null
